# Update recomputed weight/strength relation statistics (proposed & Crowd-Certain columns)
# following fixes to the enum classes used during computation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1.068436924282015
$ws.Range("E2").Value = 0.004911561765477959
$ws.Range("D3").Value = 1.041082252532498
$ws.Range("E3").Value = 0.0466368195462401
$ws.Range("D4").Value = 0.8727332866680846
$ws.Range("E4").Value = 0.1070517710685537
$ws.Range("D5").Value = 0.7226032995335685
$ws.Range("E5").Value = 0.276907900995189
$ws.Range("D6").Value = 0.8018908945232304
$ws.Range("E6").Value = 0.2280541315483521
$ws.Range("D7").Value = 0.7494966182683133
$ws.Range("E7").Value = 0.7443756448389504
$ws.Range("D8").Value = 0.8400392508243132
$ws.Range("E8").Value = 0.9591005218259062
$ws.Range("D9").Value = 0.8779887670292905
$ws.Range("E9").Value = 0.9587314673657915
$ws.Range("D10").Value = 1.046586828889832
$ws.Range("E10").Value = 1.305461351468049
$ws.Range("D11").Value = 1.07051922024726
$ws.Range("E11").Value = 1.36778317256613
$ws.Range("A12").Value = 0.6380604845384019
$ws.Range("D12").Value = 1.085817337881361
$ws.Range("E12").Value = 1.395451662880637
$ws.Range("D13").Value = 1.064392391594276
$ws.Range("E13").Value = 1.364387745422182
$ws.Range("A14").Value = 0.6503828814202761
$ws.Range("D14").Value = 1.082909529646328
$ws.Range("E14").Value = 1.385278309489174
$ws.Range("D15").Value = 1.082850435956414
$ws.Range("E15").Value = 1.386969685994254
$ws.Range("D16").Value = 1.097480809268914
$ws.Range("E16").Value = 1.405779737926417
$ws.Range("D17").Value = 1.092233264636301
$ws.Range("E17").Value = 1.403294025897283
$ws.Range("A18").Value = 0.8022805061070413
$ws.Range("D18").Value = 1.100784810245113
$ws.Range("E18").Value = 1.415069161967779
$ws.Range("D19").Value = 1.096678993760572
$ws.Range("E19").Value = 1.409387077677025
$ws.Range("D20").Value = 1.105345432149485
$ws.Range("E20").Value = 1.421359652391098
$ws.Range("D21").Value = 1.100129652062835
$ws.Range("E21").Value = 1.414008597365515
